$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8

$ws.Cells.Item($row, 1).Value = "2025-08-13 09:42:22 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-13 15:12:22 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""

# Copy formatting (style) from the previous row so the new row matches
# the look of the rest of the data rows.
$ws.Range("A7:H7").Copy()
$ws.Range("A8:H8").PasteSpecial(-4122)
$excel.CutCopyMode = 0
